$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = $null
$ws.Range("N32").Value = $null
$ws.Range("H62").Value = 1200
$ws.Range("I62").Value = 1200
$ws.Range("K62").Value = 1200
$ws.Range("M62").Value = -576
$ws.Range("H64").Value = 51098.906
$ws.Range("J64").Value = 3693.1667
$ws.Range("L64").Value = 3693.1667
$ws.Range("N64").Value = -4189.1667
$ws.Range("H65").Value = 1200
$ws.Range("I65").Value = 1200
$ws.Range("K65").Value = 6000
$ws.Range("M65").Value = -2880
$ws.Range("H67").Value = 51098.906
$ws.Range("J67").Value = 3693.1667
$ws.Range("L67").Value = 3693.1667
$ws.Range("N67").Value = -5409.1667
$ws.Range("H96").Value = 614.73334
$ws.Range("I96").Value = 368
$ws.Range("K96").Value = 1104
$ws.Range("M96").Value = 269
$ws.Range("H129").Value = 2853.66
$ws.Range("I129").Value = 7567.643
$ws.Range("J129").Value = 1020.44446
$ws.Range("K129").Value = 22702.929
$ws.Range("L129").Value = 3061.33338
$ws.Range("M129").Value = -17702.929
$ws.Range("N129").Value = -13061.33338
$ws.Range("H137").Value = 1170.2354
$ws.Range("I137").Value = 955.6047
$ws.Range("J137").Value = 2323.875
$ws.Range("K137").Value = 2866.8141
$ws.Range("L137").Value = 6971.625
$ws.Range("M137").Value = -316.8141000000001
$ws.Range("N137").Value = -12071.625
$ws.Range("H138").Value = 2123.519
$ws.Range("I138").Value = 1022.1111
$ws.Range("J138").Value = 3581.2646
$ws.Range("K138").Value = 3066.3333
$ws.Range("L138").Value = 10743.7938
$ws.Range("M138").Value = 2073.6667
$ws.Range("N138").Value = -21023.7938

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 21369.334
$ws.Range("I32").Value = 3186.4268
$ws.Range("J32").Value = 172893.56
$ws.Range("K32").Value = 3186.4268
$ws.Range("L32").Value = 172893.56
$ws.Range("M32").Value = -2899.4268
$ws.Range("N32").Value = -173467.56
$ws.Range("H122").Value = 1688.6923
$ws.Range("I122").Value = 1749.091
$ws.Range("J122").Value = 1356.5
$ws.Range("K122").Value = 5247.272999999999
$ws.Range("L122").Value = 4069.5
$ws.Range("M122").Value = -2797.272999999999
$ws.Range("N122").Value = -8969.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H103").Value = 24609.334
$ws.Range("J103").Value = 24609.334
$ws.Range("L103").Value = 24609.334
$ws.Range("N103").Value = -26953.334
$ws.Range("H107").Value = 45491904
$ws.Range("I107").Value = 58871492
$ws.Range("J107").Value = 1305.4
$ws.Range("K107").Value = 58871492
$ws.Range("L107").Value = 1305.4
$ws.Range("M107").Value = -58869572
$ws.Range("N107").Value = -5145.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1068.4
$ws.Range("I16").Value = 1052.6
$ws.Range("J16").Value = 1100
$ws.Range("K16").Value = 1052.6
$ws.Range("L16").Value = 1100
$ws.Range("M16").Value = -765.5999999999999
$ws.Range("N16").Value = -1674
$ws.Range("H62").Value = 2495
$ws.Range("I62").Value = 2027.5
$ws.Range("J62").Value = 2628.5715
$ws.Range("K62").Value = 2027.5
$ws.Range("L62").Value = 2628.5715
$ws.Range("M62").Value = -1403.5
$ws.Range("N62").Value = -3876.5715
$ws.Range("H65").Value = 2495
$ws.Range("I65").Value = 2027.5
$ws.Range("J65").Value = 2628.5715
$ws.Range("K65").Value = 10137.5
$ws.Range("L65").Value = 13142.8575
$ws.Range("M65").Value = -7017.5
$ws.Range("N65").Value = -19382.8575
$ws.Range("H94").Value = 1129.9333
$ws.Range("I94").Value = 1012
$ws.Range("K94").Value = 1012
$ws.Range("M94").Value = -561
$ws.Range("H96").Value = 20081.334
$ws.Range("J96").Value = 20081.334
$ws.Range("L96").Value = 20081.334
$ws.Range("N96").Value = -25573.334
$ws.Range("H113").Value = 1068.4
$ws.Range("I113").Value = 1052.6
$ws.Range("J113").Value = 1100
$ws.Range("K113").Value = 1052.6
$ws.Range("L113").Value = 1100
$ws.Range("M113").Value = 1117.4
$ws.Range("N113").Value = -5440

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 24000
$ws.Range("J37").Value = 24000
$ws.Range("L37").Value = 72000
$ws.Range("N37").Value = -72224
$ws.Range("H129").Value = 13164590
$ws.Range("I129").Value = 35714692
$ws.Range("J129").Value = 10364.167
$ws.Range("K129").Value = 107144076
$ws.Range("L129").Value = 31092.501
$ws.Range("M129").Value = -107139076
$ws.Range("N129").Value = -41092.501
$ws.Range("H131").Value = 823.87
$ws.Range("J131").Value = 857.90216
$ws.Range("L131").Value = 2573.70648
$ws.Range("N131").Value = -12653.70648

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 93385.87
$ws.Range("I70").Value = 139131.27
$ws.Range("K70").Value = 139131.27
$ws.Range("M70").Value = -138861.27
$ws.Range("H73").Value = 93385.87
$ws.Range("I73").Value = 139131.27
$ws.Range("K73").Value = 139131.27
$ws.Range("M73").Value = -138195.27
$ws.Range("H122").Value = 1670.64
$ws.Range("I122").Value = 1625.7273
$ws.Range("K122").Value = 4877.1819
$ws.Range("M122").Value = -2427.1819

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3194.389
$ws.Range("I7").Value = 1919.4
$ws.Range("J7").Value = 4788.125
$ws.Range("K7").Value = 1919.4
$ws.Range("L7").Value = 4788.125
$ws.Range("M7").Value = -1807.4
$ws.Range("N7").Value = -5012.125
$ws.Range("H22").Value = 1281.1904
$ws.Range("I22").Value = 639.8
$ws.Range("J22").Value = 1481.625
$ws.Range("K22").Value = 639.8
$ws.Range("L22").Value = 1481.625
$ws.Range("M22").Value = -344.8
$ws.Range("N22").Value = -2071.625
$ws.Range("H27").Value = 1281.1904
$ws.Range("I27").Value = 639.8
$ws.Range("J27").Value = 1481.625
$ws.Range("K27").Value = 639.8
$ws.Range("L27").Value = 1481.625
$ws.Range("M27").Value = -532.8
$ws.Range("N27").Value = -1695.625
$ws.Range("H40").Value = 115876.664
$ws.Range("I40").Value = 512700
$ws.Range("J40").Value = 2498.5715
$ws.Range("K40").Value = 512700
$ws.Range("L40").Value = 2498.5715
$ws.Range("M40").Value = -512564
$ws.Range("N40").Value = -2770.5715
$ws.Range("H122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").Value = $null
$ws.Range("H126").Value = 3194.389
$ws.Range("I126").Value = 1919.4
$ws.Range("J126").Value = 4788.125
$ws.Range("K126").Value = 5758.200000000001
$ws.Range("L126").Value = 14364.375
$ws.Range("M126").Value = -3288.200000000001
$ws.Range("N126").Value = -19304.375
$ws.Range("H132").Value = 3004.439
$ws.Range("I132").Value = 3211.2334
$ws.Range("J132").Value = 2440.4546
$ws.Range("K132").Value = 9633.700199999999
$ws.Range("L132").Value = 7321.3638
$ws.Range("M132").Value = -7103.700199999999
$ws.Range("N132").Value = -12381.3638

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 6143.9414
$ws.Range("J54").Value = 6912.643
$ws.Range("L54").Value = 6912.643
$ws.Range("N54").Value = -7952.643
$ws.Range("H119").Value = 42490
$ws.Range("J119").Value = 42490
$ws.Range("L119").Value = 42490
$ws.Range("N119").Value = -52166
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = $null
$ws.Range("N122").Value = $null
$ws.Range("H132").Value = 10439
$ws.Range("I132").Value = 45000
$ws.Range("K132").Value = 135000
$ws.Range("M132").Value = -132470
